$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MENT-69: Add health facility from Ile and Alto Molocue.
#
# Shared strings must be appended in the exact order the facilities /
# district were typed in the source workbook, so the new <si> entries
# line up: CS Welele, Alto Molocue, CS Mucuaba, CS Curruane, Timoteo Salua.
$ws.Range("D21").Value = "CS Welele"
$ws.Range("C21").Value = "Alto Molocué"
$ws.Range("D22").Value = "CS Mucuaba"
$ws.Range("D23").Value = "CS Curruane"
$ws.Range("E21").Value = "Timoteo Salua"

# Row 21 - Alto Molocué / CS Welele / Timoteo Salua
$ws.Range("B21").Value = 1
$ws.Range("F21").Value = 43145

# Row 22 - Ile / CS Mucuaba / Raul Ribeiro
$ws.Range("B22").Value = 5
$ws.Range("C22").Value = "Ile"
$ws.Range("E22").Value = "Raul Ribeiro"
$ws.Range("F22").Value = 43145

# Row 23 - Ile / CS Curruane / Raul Ribeiro
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = "Ile"
$ws.Range("E23").Value = "Raul Ribeiro"
$ws.Range("F23").Value = 43145

# Match the formatting used by the rest of the table (Arial body font and
# d/m/yyyy date format) and the row height of the existing data rows.
$ws.Range("B21:F23").Font.ThemeColor = 1
$ws.Range("F21:F23").NumberFormat = "d/m/yyyy"
$ws.Rows.Item(21).RowHeight = 15.75
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 15.75

# Restore the selection the author left the sheet with.
$ws.Range("D6").Select() | Out-Null
